# Apply the "prompt design results" edit:
#  1) Fix a typo in the header cell F1 ("séquencage" -> "séquençage")
#  2) Re-order the data rows (2..14) into their new positions, carrying
#     each row's full set of values (A..G) along with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Header typo fix -----------------------------------------------
$ws.Cells.Item(1, 6).Value = "Qualité du séquençage"

# --- 2) Capture the existing data rows (2..14, columns A..G) ----------
$firstDataRow = 2
$lastDataRow = 14
$lastCol = 7

$captured = @{}
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $vals = @()
    for ($c = 1; $c -le $lastCol; $c++) {
        $vals += ,$ws.Cells.Item($r, $c).Text
    }
    $captured[$r] = $vals
}

# --- 3) Mapping of new row number -> original row number ---------------
# (derived from comparing the row contents before/after the edit)
$mapping = @{
    2  = 13
    3  = 8
    4  = 5
    5  = 7
    6  = 11
    7  = 10
    8  = 2
    9  = 4
    10 = 14
    11 = 12
    12 = 6
    13 = 9
    14 = 3
}

# --- 4) Write the rows back out in their new order ----------------------
# Values that look like a plain number or a percentage (e.g. "21",
# "24219576", "10%") must be forced back to Text format before being
# re-assigned, otherwise Excel's auto-detection would silently turn them
# into numeric/percentage cells instead of keeping them as text.
$numberPattern = '^-?\d+(\.\d+)?$'
$percentPattern = '^-?\d+(\.\d+)?%$'

foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    $vals = $captured[$oldRow]
    for ($c = 1; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($newRow, $c)
        $val = $vals[$c - 1]
        if ($val -match $numberPattern -or $val -match $percentPattern) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $val
    }
}
